# Add new tsd-repo entries: PhoneGap, Sammyjs, SoundJS, sugar, toastr
# (and the pre-existing underscore-typed row), re-sorted alphabetically by name,
# into Sheet6, and refresh Sheet7's view/selection to match the new extent.

$wb = $excel.ActiveWorkbook
$sheet6 = $wb.Worksheets.Item("Sheet6")
$sheet7 = $wb.Worksheets.Item("Sheet7")

# ---------------------------------------------------------------------------
# Final (sorted) data for rows 38-49 of Sheet6.
# Columns: A name, B description, C version, D key, E author, F url
# ---------------------------------------------------------------------------
$rows = @(
    @{ r = 38; A = "PhoneGap";         B = "Framework that supports 7 mobile platforms"; C = "2.2";  D = "53535BA3-1A74-4DA2-9B98-46DFEFA74BB2"; E = "Boris Yankov";       F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/phonegap-2.2.d.ts" },
    @{ r = 39; A = "qunit";             B = "JavaScript runtime for build applications.";  C = "1.10"; D = "BD09CCAE-87F7-49FF-9CFE-5FB9F95369AB"; E = "Diullei Gomes";      F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/qunit-1.10.d.ts" },
    @{ r = 40; A = "raphael";           B = "A small JavaScript library that should simplify your work with vector graphics on the web."; C = "2.1"; D = "E6A34E99-68E8-42F8-9479-ABA9C575319C"; E = "CheCoxshall"; F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/raphael-2.1.d.ts" },
    @{ r = 41; A = "Sammyjs";           B = "RESTful Evented JavaScript"; C = "0.7"; D = "3FBC8337-4741-416D-B576-B55221CDA93A"; E = "Boris Yankov"; F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/sammyjs-0.7.d.ts" },
    @{ r = 42; A = "SoundJS";           B = "Javascript library that provides a simple API, and powerful features to make working with audio a breeze."; C = "3.0"; D = "E6459F1B-6722-4B84-9E42-150F0D7740B1"; E = "Pedro Ferreira"; F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/soundjs-0.3.d.ts" },
    @{ r = 43; A = "spin";              B = "A spinning activity indicator."; C = "1.2"; D = "3844D795-5DAC-4AE2-9B3D-37DD342E2409"; E = "Boris Yankov"; F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/spin-1.2.d.ts" },
    @{ r = 44; A = "sugar";             B = "Sugar is a Javascript library that extends native objects with helpful methods."; C = "1.3"; D = "3479CE5E-C7C5-446B-9849-36DC418DC838"; E = "Josh Baldwin"; F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/sugar-1.3.d.ts" },
    @{ r = 45; A = "TeeChart";          B = "TeeChart Pro provides complete, quick and easy to use charting"; C = "1.3"; D = "CD94F66C-57A3-4FCA-B8BF-EA93D9D1B843"; E = "Steema Software SL."; F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/teechart-1.3.d.ts" },
    @{ r = 46; A = "toastr";            B = "toastr is a Javascript library for Gnome / Growl type non-blocking notifications."; C = "1.0"; D = "B4D5487F-E37D-4326-9E08-A63B745FEFC5"; E = "Boris Yankov"; F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/toastr-1.0.d.ts" },
    @{ r = 47; A = "TweenJS";           B = "A Javascript library for tweening and animating HTML5 and Javascript properties."; C = "0.3"; D = "43E6CDBC-EAFC-4F54-A57E-C2CF2F287624"; E = "Pedro Ferreira"; F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/tweenjs-0.3.d.ts" },
    @{ r = 48; A = "underscore";        B = "Utility-belt library for JavaScript."; C = "1.4"; D = "61302984-20B9-412B-9FB3-FDFC7554E144"; E = "Boris Yankov"; F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/underscore-1.4.d.ts" },
    @{ r = 49; A = "underscore-typed";  B = "Utility-belt library for JavaScript."; C = "1.4"; D = "9C68EE32-8B0A-401E-B35F-C2F7C0761A7C"; E = "Josh Baldwin"; F = "https://github.com/borisyankov/DefinitelyTyped/raw/master/Definitions/underscore-typed-1.4.d.ts" }
)

foreach ($row in $rows) {
    $r = $row.r
    $sheet6.Cells.Item($r, 1).Value = $row.A
    $sheet6.Cells.Item($r, 2).Value = $row.B
    $sheet6.Cells.Item($r, 3).Value = $row.C
    $sheet6.Cells.Item($r, 4).Value = $row.D
    $sheet6.Cells.Item($r, 5).Value = $row.E
    $sheet6.Cells.Item($r, 6).Value = $row.F
}

# Highlight the genuinely new rows / cells the same way the author did
# (existing "new entry" style index 2 used elsewhere in the sheet).
$sheet6.Range("A38").Style = $sheet6.Range("A41").Style
$sheet6.Range("A41").Style = $sheet6.Range("A45").Style
$sheet6.Range("A42").Style = $sheet6.Range("A45").Style
$sheet6.Range("E42").Style = $sheet6.Range("A45").Style
$sheet6.Range("A44").Style = $sheet6.Range("A45").Style
$sheet6.Range("E44").Style = $sheet6.Range("A45").Style
$sheet6.Range("A46").Style = $sheet6.Range("A45").Style
$sheet6.Range("E49").Style = $sheet6.Range("A45").Style

# ---------------------------------------------------------------------------
# View/selection bookkeeping to match the author's final workbook state:
# Sheet6 becomes the active/selected tab with A2 selected, Sheet7's selection
# widens to the new data extent.
# ---------------------------------------------------------------------------
$sheet6.Activate()
$sheet6.Range("A2").Select()

$sheet7.Activate()
$sheet7.Range("B2:B49").Select()
$sheet6.Activate()
